$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header B1 from "U_S_census" to "US_census" (C1 "UN" stays as-is)
$ws.Range("B1").Value = "US_census"

# Columns B and C currently hold raw population counts; convert them to
# billions (divide by 1e9), row by row, for all existing data rows.
for ($r = 2; $r -le 67; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $bVal = $bCell.Value2
    if ($bVal -ne $null) {
        $bCell.Value = $bVal / 1000000000
    }

    $cCell = $ws.Cells.Item($r, 3)
    $cVal = $cCell.Value2
    if ($cVal -ne $null) {
        $cCell.Value = $cVal / 1000000000
    }
}

# Row 68 only has a U.S. census figure; scale it and add a new UN figure of 0.
$b68 = $ws.Cells.Item(68, 2)
$b68.Value = $b68.Value2 / 1000000000
$ws.Cells.Item(68, 3).Value = 0

# Update the active selection to match the saved view.
$ws.Range("I12").Select()
